# changed kill to eliminate; changed stick to what is left
$wb = $excel.ActiveWorkbook

# --- gameData: "play" flag (C2) flips on now that kill -> eliminate switch is live ---
$gameSheet = $wb.Worksheets.Item("gameData")
$gameSheet.Range("C2").Value = $true

# --- userTrackerData: populate with the current tracked session/user row ---
$trackSheet = $wb.Worksheets.Item("userTrackerData")
$headerSource = $wb.Worksheets.Item("gameData").Range("B1")

# Header row (bold / bordered / centered style, matching the other data sheets)
$trackSheet.Range("B1").Value = "username"
$trackSheet.Range("C1").Value = "state"
$trackSheet.Range("D1").Value = "db"
$trackSheet.Range("E1").Value = "chat_id"
$trackSheet.Range("F1").Value = "elimination_target"

$headerSource.Copy()
$trackSheet.Range("B1:F1").PasteSpecial(-4122)

# Data row
$trackSheet.Range("A2").Value = "praveeeenk"
$headerSource.Copy()
$trackSheet.Range("A2").PasteSpecial(-4122)

$trackSheet.Range("B2").Value = "praveeeenk"
$trackSheet.Range("C2").Value = ""
$trackSheet.Range("D2").Value = "<dbhelper.DBHelper object at 0x00000156784073A0>"
$trackSheet.Range("E2").Value = 258884638
$trackSheet.Range("F2").Value = ""

$excel.CutCopyMode = $false
